# Apply the "rebrand Alvearie -> LinuxForHealth" edit described by the diff.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-allowed-amount-outpatient"

# Version
$meta.Range("B3").Value = "8.0.0"

# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 (root "Extension" element) no longer carries the ele-1/ext-1
# constraint text in its Constraint(s) column.
$elements.Range("AI2").Value = ""

# Row 5 ("Extension.url" Fixed Value) reflects the same URL rebrand.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-allowed-amount-outpatient"
